$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "77.887.86"
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").Value = "3.167.86"
$ws.Range("E3").Value = "  +6.64%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "203.57"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "630.37"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.229"
$ws.Range("E8").Value = "  +14.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.572"
$ws.Range("E9").Value = "  +4.61%  "
$ws.Range("D10").Value = "3.167.43"
$ws.Range("E10").Value = "  +6.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.566"
$ws.Range("E11").Value = "  +31.47%  "
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.42"
$ws.Range("E13").Value = "  +9.07%  "
$ws.Range("D14").Value = "3.750.02"
$ws.Range("E14").Value = "  +6.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000226"
$ws.Range("E15").Value = "  +20.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "31.55"
$ws.Range("E16").Value = "  +8.76%  "
$ws.Range("D17").Value = "77.790.96"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "3.163.76"
$ws.Range("E18").Value = "  +7.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.29"
$ws.Range("E19").Value = "  +6.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.36"
$ws.Range("E20").Value = "  +6.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "428.90"
$ws.Range("E21").Value = "  +15.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.84"
$ws.Range("E22").Value = "  +25.68%  "
$ws.Range("E23").Value = "  +13.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.71"
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("D25").Value = "3.324.41"
$ws.Range("E25").Value = "  +6.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.71"
$ws.Range("E26").Value = "  +9.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "76.29"
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.92"
$ws.Range("E28").Value = "  +12.57%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +8.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.84"
$ws.Range("E32").Value = "  +7.55%  "
$ws.Range("E33").Value = "  +7.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "519.97"
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("E35").Value = "  +2.50%  "
$ws.Range("E36").Value = "  +23.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "22.37"
$ws.Range("E37").Value = "  +10.24%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.396"
$ws.Range("E39").Value = "  +4.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.61"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "196.03"
$ws.Range("E41").Value = "  +6.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.06"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.40"
$ws.Range("E45").Value = "  +9.91%  "
$ws.Range("E46").Value = "  +13.49%  "
$ws.Range("E47").Value = "  +7.96%  "
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "42.89"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  +10.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.622"
$ws.Range("E51").Value = "  +6.48%  "
